$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 6 data)
$ws.Range("D2").Value = 44519
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1240
$ws.Range("N2").Value = '$/kilo'
$ws.Range("O2").Value = "Provincia de Linares"
$ws.Range("P2").Value = 1240

# Row 3 (was row 5 data)
$ws.Range("D3").Value = 44468
$ws.Range("H3").Value = "Verde"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1920
$ws.Range("N3").Value = '$/kilo'
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1920

# Row 4 (was row 9 data)
$ws.Range("D4").Value = 44545
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 1700
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1755
$ws.Range("N4").Value = '$/kilo'
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 1755

# Row 5 (was row 4 data)
$ws.Range("D5").Value = 44526
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1600
$ws.Range("M5").Value = 1550
$ws.Range("N5").Value = '$/kilo'
$ws.Range("O5").Value = "Provincia de Linares"
$ws.Range("P5").Value = 1550

# Row 6 (was row 11 data)
$ws.Range("D6").Value = 44510
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = 1350
$ws.Range("N6").Value = '$/kilo'
$ws.Range("O6").Value = "Provincia de Linares"
$ws.Range("P6").Value = 1350

# Row 7 (was row 2 data)
$ws.Range("D7").Value = 44477
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1460
$ws.Range("N7").Value = '$/kilo'
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 1460

# Row 8 (was row 3 data)
$ws.Range("D8").Value = 44524
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1600
$ws.Range("M8").Value = 1550
$ws.Range("N8").Value = '$/kilo'
$ws.Range("O8").Value = "Provincia de Talca"
$ws.Range("P8").Value = 1550

# Row 9 (was row 8 data)
$ws.Range("D9").Value = 44511
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 1300
$ws.Range("L9").Value = 1400
$ws.Range("M9").Value = 1350
$ws.Range("N9").Value = '$/kilo'
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1350

# Row 10 (was row 7 data)
$ws.Range("D10").Value = 44496
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 550
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1773
$ws.Range("N10").Value = '$/paquete'
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1773

# Row 11 (was row 10 data)
$ws.Range("D11").Value = 44489
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 1400
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1450
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = "Provincia de Linares"
$ws.Range("P11").Value = 1450
